$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.487815260887146
$ws.Range("B1").Value = 3.661987066268921
$ws.Range("C1").Value = 2.20524263381958
$ws.Range("D1").Value = 1.274955868721008
$ws.Range("E1").Value = 0.7593414187431335
